# Word COM-interop script implementing the tracked changes from the diff.
$d = $word.ActiveDocument

# wdReplaceAll = 2, wdFindContinue = 1
$wdReplaceAll = 2
$wdFindContinue = 1

# --- 1. "...480,96 руб., время..." -> "...480,96 рублей, время..." ---
$d.Content.Find.Execute(
    "Минимальная заработная плата с 01.04.2022 г. составляет 480,96 руб., время на разработку программного обеспечения 208 часов (итог по графе 4 таблицы 5.1",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Минимальная заработная плата с 01.04.2022 г. составляет 480,96 рублей, время на разработку программного обеспечения 208 часов (итог по графе 4 таблицы 5.1",
    $wdReplaceAll) | Out-Null

# --- 2. "...595,47 руб." -> "...595,47 рублей." ---
$d.Content.Find.Execute(
    "Таким образом, заработная плата составит 595,47 руб.",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Таким образом, заработная плата составит 595,47 рублей.",
    $wdReplaceAll) | Out-Null

# --- 3. "руб [" -> "рублей [" ---
$d.Content.Find.Execute(
    "руб [",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "рублей [",
    $wdReplaceAll) | Out-Null

# --- 4. ") составляет 0,40220 руб./кВтч " -> ") составляет 0,40220 руб/кВтч " ---
$d.Content.Find.Execute(
    ") составляет 0,40220 руб./кВтч ",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    ") составляет 0,40220 руб/кВтч ",
    $wdReplaceAll) | Out-Null

# --- 5. "1 руб." -> "1 рублей." ---
$d.Content.Find.Execute(
    "1 руб.",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "1 рублей.",
    $wdReplaceAll) | Out-Null

# --- 6. " руб. " -> " рублей." ---
$d.Content.Find.Execute(
    " руб. ",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    " рублей.",
    $wdReplaceAll) | Out-Null

# --- 7. bare "руб." (the one after "...1839,77 ") -> "рублей." ---
# Scope the search to the one paragraph that ends in a lone "руб." run
# (there are several other "руб."-ending paragraphs that must stay as-is
# because they are not touched by this edit).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*1839,77*") {
        $p.Range.Find.Execute(
            "руб.",
            $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
            "рублей.",
            $wdReplaceAll) | Out-Null
    }
}

Write-Host "Text edits done"
